$d = $word.ActiveDocument

# 1. Remove the "Berichtsdatum: 22. Januar 2024" paragraph together with the
#    following blank (single-space) paragraph, including both paragraph marks,
#    so the remaining "In den letzten Wochen..." paragraph keeps its own
#    (non-bold) paragraph properties.
$datePara = $null
$spacePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Berichtsdatum: 22. Januar 2024") {
        $datePara = $p
        $spacePara = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($datePara -ne $null -and $spacePara -ne $null) {
    $rng = $d.Range($datePara.Range.Start, $spacePara.Range.End)
    $rng.Delete()
}

# 2. Update the closing sentence with the revised wording.
$old = "Wenn das Produkt sein derzeitiges Verkaufstempo in der Welt der Gesundheits- und Fitnessbranche beibehalten kann, ist es vielleicht reif für eine nationale Vermarktung."
$new = "Wenn das Produkt sein derzeitiges Verkaufstempo in der Hochburg der Gesundheits- und Fitnesswelt beibehalten kann, könnte es für eine landesweite Vermarktung bereit sein."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
